$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.685.77'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.63%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.689.64'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.62%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.47%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.27%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.71%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3702'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.54%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '48.52'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.81%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3383'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.86%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.193'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.83%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07387'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.54%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.78%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.233'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.22%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.66'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.62%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.890'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.03%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.688.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.59%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001109'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.87%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06677'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.71%  '

# Row 19
$ws.Range('E19').Value = '  +0.64%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '82.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.20%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.44%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.254'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.15%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.68%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.691.10'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.81%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.444'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.74%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.724'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.87%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.00%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '147.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.17%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '130.95'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.05%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.875.68'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.64%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.224'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +25.61%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.593'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.86%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.190'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.39%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '13.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.21%  '

# Row 35
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08634'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.81%  '

# Row 36
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.738'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.81%  '

# Row 37
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06537'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.79%  '

# Row 38
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.446'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.89%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.918'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.96%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02372'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.97%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2183'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.87%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.247'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.26%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6318'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.06%  '

# Row 44
$ws.Range('E44').Value = '  +0.67%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.19%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.795'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.05%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5983'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.29%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.073'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.63%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '127.19'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.70%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07195'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.13%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.39'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.77%  '
